$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the header parameters (commit: "Changed name of parameter.")
$ws.Range("C1").Value = "begin_part"
$ws.Range("D1").Value = "end_part"

# Update the selected cell to match the new active cell/selection
$ws.Range("D11").Select()
